$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Req 1 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "The Solution must support missiles launched by the user"
$ws.Range("C2").Value = "Mandatory"
$ws.Range("D2").Value = "Functional"
$ws.Range("E2").Value = "Gameplay"
$ws.Range("F2").Value = "Diagram 1, main loop detects key stroke, Diagram 2 draws and spawns missiles."

# --- Sub items 1.a - 1.d (A3:A6) ---
$ws.Range("A3").Value = "1.a"
$ws.Range("A4").Value = "1.b"
$ws.Range("A5").Value = "1.c"
$ws.Range("A6").Value = "1.d"

# --- Sub items 1.e - 1.k (A7:A13) formatted as Text ---
$ws.Range("A7:A13").NumberFormat = "@"
$ws.Range("A7").Value = "1.e"
$ws.Range("A8").Value = "1.f"
$ws.Range("A9").Value = "1.g"
$ws.Range("A10").Value = "1.h"
$ws.Range("A11").Value = "1.i"
$ws.Range("A12").Value = "1.j"
$ws.Range("A13").Value = "1.k"

# --- Row 14 : Req 2 ---
$ws.Range("A14").Value = 2
$ws.Range("A14").NumberFormat = "@"
$ws.Range("B14").Value = "The Solution must support Ships launched by the solution"
$ws.Range("C14").Value = "Mandatory"
$ws.Range("D14").Value = "Functional"
$ws.Range("E14").Value = "Gameplay"

# --- Sub items 2.a - 2.d.iii.4 (A15:A25) formatted as Text ---
$ws.Range("A15:A25").NumberFormat = "@"
$ws.Range("A15").Value = "2.a"
$ws.Range("A16").Value = "2.b"
$ws.Range("A17").Value = "2.c"
$ws.Range("A18").Value = "2.d"
$ws.Range("A19").Value = "2.d.i"
$ws.Range("A20").Value = "2.d.ii"
$ws.Range("A21").Value = "2.d.iii"
$ws.Range("A22").Value = "2.d.iii.1"
$ws.Range("A23").Value = "2.d.iii.2"
$ws.Range("A24").Value = "2.d.iii.3"
$ws.Range("A25").Value = "2.d.iii.4"

# --- Rows 26-29 : Req 3-6 ---
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "The Solution must end the game when hit count reaches ten"
$ws.Range("C26").Value = "Mandatory"
$ws.Range("D26").Value = "Functional"
$ws.Range("E26").Value = "End game"

$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "The Solution must end the game when user clicks left mouse button"
$ws.Range("C27").Value = "Mandatory"
$ws.Range("D27").Value = "Functional"
$ws.Range("E27").Value = "End game"

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "The Solution must end the game when user clicks esc key"
$ws.Range("C28").Value = "Mandatory"
$ws.Range("D28").Value = "Functional"
$ws.Range("E28").Value = "End game"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "The Solution must end the game if user hasn't launched a missile in last 5 minutes"
$ws.Range("C29").Value = "Mandatory"
$ws.Range("D29").Value = "Functional"
$ws.Range("E29").Value = "End game"

$ws.Range("A26:A29").NumberFormat = "@"

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 71.77734375

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("H5").Select()
